# Auto-generated edit script
$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in @(1, 4)) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # 1. Insert a new row at position 18 (shifts old rows 18-47 down to 19-48)
    $ws.Rows.Item(18).Insert()

    # 2. Populate the newly inserted row 18 with the new event data
    $ws.Range("A18").Value = 17
    $ws.Range("B18").Value = '2024-02-12'
    $ws.Range("C18").Value = '赣州·宅舞联萌·随舞动漫派对（免费活动)'
    $ws.Range("D18").Value = '金岭东大道新都汇西侧约100米 万达广场'
    $ws.Range("E18").Value = '2024.02.12 14:00-02.13 19:00'
    $ws.Range("F18").Value = 0
    $ws.Range("G18").Value = 22.33
    $ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=81540'
    $ws.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202401/5gLDYtbv1706608938962.jpeg'
    $ws.Range("A18").Style = $ws.Range("A19").Style

    # 3. Apply F-column (want-to-go count) corrections for rows shifted down (post-insert row numbers)
    $ws.Range("F19").Value = 47
    $ws.Range("F20").Value = 46
    $ws.Range("F21").Value = 92
    $ws.Range("F22").Value = 897
    $ws.Range("F23").Value = 1391
    $ws.Range("F26").Value = 184
    $ws.Range("F27").Value = 70
    $ws.Range("F29").Value = 38
    $ws.Range("F31").Value = 210
    $ws.Range("F33").Value = 267
    $ws.Range("F34").Value = 1608
    $ws.Range("F38").Value = 575
    $ws.Range("F40").Value = 3573
    $ws.Range("F42").Value = 188
    $ws.Range("F43").Value = 902
    $ws.Range("F46").Value = 63
    $ws.Range("F47").Value = 42

    # 4. Apply F-column corrections for rows 2-17 (not affected by the insert)
    $ws.Range("F7").Value = 79
    $ws.Range("F8").Value = 451
    $ws.Range("F11").Value = 564
    $ws.Range("F15").Value = 363
    $ws.Range("F17").Value = 90

    # 5. Re-normalize column A (sequence number = row - 1) for every data row, 1..48
    for ($r = 1; $r -le 48; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}
